# Generate Report for Handback
# Refresh the handback-status report for the e9de8414-... source file: it has
# been re-handed-off and handed back for both the zh-cn and de-de targets, so
# stamp the new handoff/handback timestamps on each locale sheet and bump the
# "Latest HO Xliff Generate Date" on the Overview sheet to the newest of the
# two new handoff timestamps (de-de's 12:50:06).

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-22 12:50:06"

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-08-22 12:49:56"
$zhcn.Range("K3").Value = "2016-08-22 12:50:34"

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-08-22 12:50:06"
$dede.Range("K3").Value = "2016-08-22 12:50:41"
